$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values (Motor / Chasis / Patente reference numbers and plazo)
$ws.Range("W2").Value = "MMM115"
$ws.Range("X2").Value = "MASDAS12316"
$ws.Range("Y2").Value = "ASDAKE1236"
$ws.Range("H2").Value = "Anual"
$ws.Range("E2").Value = 8684079401

# Adjust column R width slightly
$ws.Columns.Item(18).ColumnWidth = 5.83

# Update the active selection / view
$ws.Range("B7").Select()
